$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-05 02:28:10"

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
